# Update the "Förändrad" date column (C) for rows 2-7 from 2023-09-06 (45175)
# to 2023-09-14 (45183), keeping the existing cell formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45183
    }
}
